# Weekly refresh of the Alcachofa (Hortaliza, Terminal Hortofrutícola Agro
# Chillán) price rows: dates and their associated volume/price/origin data
# are rotated to the following week's reported figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44432
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("P2").Value = 362
$ws.Range("D3").Value = 44484
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("P3").Value = 288
$ws.Range("D4").Value = 44446
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 12500
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12750
$ws.Range("P4").Value = 319
$ws.Range("D5").Value = 44425
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 362
$ws.Range("D6").Value = 44420
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 13000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 13500
$ws.Range("P6").Value = 338
$ws.Range("D7").Value = 44468
$ws.Range("H7").Value = "Madrigal"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 12500
$ws.Range("P7").Value = 312
$ws.Range("D8").Value = 44399
$ws.Range("H8").Value = "Española"
$ws.Range("I8").Value = "Segunda"
$ws.Range("K8").Value = 15500
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15750
$ws.Range("P8").Value = 394
$ws.Range("D9").Value = 44487
$ws.Range("J9").Value = 100
$ws.Range("D10").Value = 44455
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13500
$ws.Range("P10").Value = 338
$ws.Range("D11").Value = 44435
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("P11").Value = 362
$ws.Range("D12").Value = 44516
$ws.Range("K12").Value = 11000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 11500
$ws.Range("O12").Value = "Provincia del Elquí"
$ws.Range("P12").Value = 288
$ws.Range("D13").Value = 44508
$ws.Range("J13").Value = 160
$ws.Range("D14").Value = 44417
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 16000
$ws.Range("M14").Value = 15500
$ws.Range("P14").Value = 388
$ws.Range("D15").Value = 44453
$ws.Range("J15").Value = 160
$ws.Range("K15").Value = 12500
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12750
$ws.Range("P15").Value = 319
$ws.Range("D16").Value = 44475
$ws.Range("J16").Value = 120
$ws.Range("D17").Value = 44495
$ws.Range("D18").Value = 44505
$ws.Range("D19").Value = 44510
$ws.Range("K19").Value = 11000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 11500
$ws.Range("P19").Value = 288
$ws.Range("D20").Value = 44498
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 10500
$ws.Range("L20").Value = 11000
$ws.Range("M20").Value = 10750
$ws.Range("P20").Value = 269
$ws.Range("D21").Value = 44467
$ws.Range("J21").Value = 160
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("D22").Value = 44489
$ws.Range("J22").Value = 120
$ws.Range("D23").Value = 44496
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 11000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 11500
$ws.Range("P23").Value = 288
$ws.Range("D24").Value = 44515
$ws.Range("J24").Value = 120
$ws.Range("K24").Value = 11000
$ws.Range("L24").Value = 12000
$ws.Range("M24").Value = 11500
$ws.Range("P24").Value = 288
$ws.Range("D25").Value = 44426
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 14000
$ws.Range("M25").Value = 13500
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 338
$ws.Range("D26").Value = 44488
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 11000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 11500
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 288
$ws.Range("D27").Value = 44473
$ws.Range("O27").Value = "Provincia del Elquí"
$ws.Range("D28").Value = 44482
$ws.Range("J28").Value = 120
$ws.Range("D29").Value = 44454
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 13000
$ws.Range("L29").Value = 14000
$ws.Range("M29").Value = 13500
$ws.Range("P29").Value = 338
$ws.Range("D30").Value = 44490
$ws.Range("J30").Value = 100
$ws.Range("D31").Value = 44427
$ws.Range("D32").Value = 44491
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 11000
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = 11500
$ws.Range("P32").Value = 288
$ws.Range("D33").Value = 44494
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 11500
$ws.Range("P33").Value = 288
$ws.Range("D34").Value = 44503
$ws.Range("D35").Value = 44512
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 11500
$ws.Range("P35").Value = 288
